# Slide 7, shape 2 ("Tijdelijke aanduiding voor inhoud 2") contains the
# command-line snippet that used to read "$ conda activate science" split
# across three runs ("$ " / "conda" / " activate science"). The author
# collapsed it into a single run reading "$ source activate science".

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(7)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

# Locate the paragraph that still has the old wording (Paragraphs(i).Text
# carries a trailing paragraph-mark character, so match on substring).
$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $tr.Paragraphs($i)
    if ($candidate.Text.Contains("conda activate science")) {
        $para = $candidate
        break
    }
}

# Keep the formatting of the paragraph's first run (the "$ " run) and fold
# the whole paragraph into that single run, dropping the other runs so the
# paragraph ends up as one run of text, matching the target markup.
$firstRun = $tr.Characters($para.Start, 2)
$rest     = $tr.Characters($para.Start + 2, $para.Length - 2)
$rest.Text = ""
$firstRun.Text = "$ source activate science"
